$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the student's name on the "Name:" line
$ws.Range("A2").Value = "Name: Abraham Khan"

# Fill in earned points (column D) for each rubric row, matching the
# possible points in column C except for row 10, which was marked 0.
$ws.Range("D4").Value = 20
$ws.Range("D5").Value = 10
$ws.Range("D6").Value = 5
$ws.Range("D7").Value = 15
$ws.Range("D8").Value = 10
$ws.Range("D9").Value = 10
$ws.Range("D10").Value = 0
$ws.Range("D11").Value = 6
$ws.Range("D12").Value = 6
$ws.Range("D13").Value = 10

# Update the view: scroll down a bit and move the active selection
$excel.Goto($ws.Range("A4"), $true)
$ws.Range("J12").Select()

$wb.Application.Calculate()
